# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = "27.462.43"
$dCell.Style = "Normal"
$ws.Range("E2").Value = "  -2.54%  "

$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = "1.747.31"
$dCell.Style = "Normal"
$ws.Range("E3").Value = "  -2.95%  "

$dCell = $ws.Range("D4")
$dCell.NumberFormat = "@"
$dCell.Value = "1.003"
$dCell.Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "

$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = "323.83"
$dCell.Style = "Normal"
$ws.Range("E5").Value = "  -0.15%  "

$ws.Range("E6").Value = "  +0.06%  "

$dCell = $ws.Range("D7")
$dCell.NumberFormat = "@"
$dCell.Value = "0.4424"
$dCell.Style = "Normal"
$ws.Range("E7").Value = "  +3.04%  "

$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = "0.3603"
$dCell.Style = "Normal"
$ws.Range("E8").Value = "  -0.86%  "

$dCell = $ws.Range("D9")
$dCell.NumberFormat = "@"
$dCell.Value = "0.07448"
$dCell.Style = "Normal"
$ws.Range("E9").Value = "  -1.59%  "

$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = "42.06"
$dCell.Style = "Normal"
$ws.Range("E10").Value = "  -6.04%  "

$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = "1.094"
$dCell.Style = "Normal"
$ws.Range("E11").Value = "  -2.58%  "

$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value = "1.002"
$dCell.Style = "Normal"
$ws.Range("E12").Value = "  +0.12%  "

$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = "20.49"
$dCell.Style = "Normal"
$ws.Range("E13").Value = "  -5.76%  "

$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value = "5.999"
$dCell.Style = "Normal"
$ws.Range("E14").Value = "  -3.33%  "

$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = "7.113"
$dCell.Style = "Normal"
$ws.Range("E15").Value = "  -3.56%  "

$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = "1.752.54"
$dCell.Style = "Normal"
$ws.Range("E16").Value = "  -3.63%  "

$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = "91.80"
$dCell.Style = "Normal"
$ws.Range("E17").Value = "  -1.22%  "

$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = "0.00001057"
$dCell.Style = "Normal"
$ws.Range("E18").Value = "  -1.22%  "

$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = "0.06400"
$dCell.Style = "Normal"
$ws.Range("E19").Value = "  +0.60%  "

$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = "1.002"
$dCell.Style = "Normal"
$ws.Range("E20").Value = "  +0.04%  "

$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = "16.75"
$dCell.Style = "Normal"
$ws.Range("E21").Value = "  -3.06%  "

$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = "5.851"
$dCell.Style = "Normal"
$ws.Range("E22").Value = "  -2.56%  "

$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = "27.518.46"
$dCell.Style = "Normal"
$ws.Range("E23").Value = "  -2.39%  "

$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = "11.14"
$dCell.Style = "Normal"
$ws.Range("E24").Value = "  -2.54%  "

$ws.Range("E25").Value = "  -2.75%  "

$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value = "161.47"
$dCell.Style = "Normal"
$ws.Range("E26").Value = "  +0.91%  "

$dCell = $ws.Range("D27")
$dCell.NumberFormat = "@"
$dCell.Value = "20.34"
$dCell.Style = "Normal"
$ws.Range("E27").Value = "  -0.47%  "

$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = "1.952.84"
$dCell.Style = "Normal"
$ws.Range("E28").Value = "  -3.38%  "

$dCell = $ws.Range("D29")
$dCell.NumberFormat = "@"
$dCell.Value = "2.074"
$dCell.Style = "Normal"
$ws.Range("E29").Value = "  -7.51%  "

$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = "124.00"
$dCell.Style = "Normal"
$ws.Range("E30").Value = "  -3.29%  "

$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = "1.067"
$dCell.Style = "Normal"
$ws.Range("E31").Value = "  -9.60%  "

$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = "3.650"
$dCell.Style = "Normal"
$ws.Range("E32").Value = "  +3.43%  "

$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = "0.08962"
$dCell.Style = "Normal"
$ws.Range("E33").Value = "  -0.80%  "

$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = "5.480"
$dCell.Style = "Normal"
$ws.Range("E34").Value = "  -7.04%  "

$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = "11.91"
$dCell.Style = "Normal"
$ws.Range("E35").Value = "  -7.26%  "

$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = "0.02293"
$dCell.Style = "Normal"
$ws.Range("E36").Value = "  -2.89%  "

$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"
$dCell.Value = "0.2075"
$dCell.Style = "Normal"
$ws.Range("E37").Value = "  -2.53%  "

$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = "0.6315"
$dCell.Style = "Normal"
$ws.Range("E38").Value = "  -3.08%  "

$dCell = $ws.Range("D39")
$dCell.NumberFormat = "@"
$dCell.Value = "0.05955"
$dCell.Style = "Normal"
$ws.Range("E39").Value = "  -2.80%  "

$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = "4.892"
$dCell.Style = "Normal"
$ws.Range("E40").Value = "  -4.92%  "

$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = "1.200"
$dCell.Style = "Normal"
$ws.Range("E41").Value = "  +0.45%  "

$ws.Range("E42").Value = "  +0.02%  "

$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = "1.385"
$dCell.Style = "Normal"
$ws.Range("E43").Value = "  -3.16%  "

$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = "7.737"
$dCell.Style = "Normal"
$ws.Range("E44").Value = "  -2.95%  "

$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = "13.25"
$dCell.Style = "Normal"
$ws.Range("E45").Value = "  -2.24%  "

$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = "3.710"
$dCell.Style = "Normal"
$ws.Range("E46").Value = "  -0.01%  "

$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = "0.5860"
$dCell.Style = "Normal"
$ws.Range("E47").Value = "  -2.77%  "

$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = "120.84"
$dCell.Style = "Normal"
$ws.Range("E48").Value = "  -3.82%  "

$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = "1.938"
$dCell.Style = "Normal"
$ws.Range("E49").Value = "  -2.71%  "

$ws.Range("E50").Value = "  -1.15%  "

$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = "0.06858"
$dCell.Style = "Normal"
$ws.Range("E51").Value = "  -1.69%  "
